# Auto-generated edit script: updates market-price derived columns (H-N)
# on multiple sheets to match the scheduled-runner data refresh described in the commit.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(21, 9).Value = 4500
$ws.Cells.Item(21, 10).Value = 2000
$ws.Cells.Item(21, 11).Value = 4500
$ws.Cells.Item(21, 12).Value = 2000
$ws.Cells.Item(21, 13).Value = -4032
$ws.Cells.Item(21, 14).Value = -2936
$ws.Cells.Item(23, 9).Value = 4500
$ws.Cells.Item(23, 10).Value = 2000
$ws.Cells.Item(23, 11).Value = 4500
$ws.Cells.Item(23, 12).Value = 2000
$ws.Cells.Item(23, 13).Value = -4266
$ws.Cells.Item(23, 14).Value = -2468
$ws.Cells.Item(32, 8).Value = 5474.25
$ws.Cells.Item(32, 9).Value = 6900
$ws.Cells.Item(32, 10).Value = 4999
$ws.Cells.Item(32, 11).Value = 6900
$ws.Cells.Item(32, 12).Value = 4999
$ws.Cells.Item(32, 13).Value = -6574
$ws.Cells.Item(32, 14).Value = -5651
$ws.Cells.Item(40, 8).Value = 2376.0715
$ws.Cells.Item(40, 9).Value = 2649.8
$ws.Cells.Item(40, 11).Value = 2649.8
$ws.Cells.Item(40, 13).Value = -2474.8
$ws.Cells.Item(64, 8).Value = 0
$ws.Cells.Item(64, 10).Value = 0
$ws.Cells.Item(64, 12).Value = 0
$ws.Cells.Item(64, 14).ClearContents()
$ws.Cells.Item(67, 8).Value = 0
$ws.Cells.Item(67, 10).Value = 0
$ws.Cells.Item(67, 12).Value = 0
$ws.Cells.Item(67, 14).ClearContents()
$ws.Cells.Item(80, 8).Value = 2871.5
$ws.Cells.Item(80, 9).Value = 1434.4
$ws.Cells.Item(80, 10).Value = 5266.6665
$ws.Cells.Item(80, 11).Value = 4303.200000000001
$ws.Cells.Item(80, 12).Value = 15799.9995
$ws.Cells.Item(80, 13).Value = -3305.200000000001
$ws.Cells.Item(80, 14).Value = -17795.9995
$ws.Cells.Item(83, 8).Value = 2871.5
$ws.Cells.Item(83, 9).Value = 1434.4
$ws.Cells.Item(83, 10).Value = 5266.6665
$ws.Cells.Item(83, 11).Value = 12909.6
$ws.Cells.Item(83, 12).Value = 47399.9985
$ws.Cells.Item(83, 13).Value = -7917.6
$ws.Cells.Item(83, 14).Value = -57383.9985
$ws.Cells.Item(97, 8).Value = 2500
$ws.Cells.Item(97, 10).Value = 2500
$ws.Cells.Item(97, 12).Value = 7500
$ws.Cells.Item(97, 14).Value = -8492
$ws.Cells.Item(116, 8).Value = 7301.6665
$ws.Cells.Item(116, 9).Value = 3905
$ws.Cells.Item(116, 11).Value = 3905
$ws.Cells.Item(116, 13).Value = -463
$ws.Cells.Item(138, 8).Value = 4765123.5
$ws.Cells.Item(138, 9).Value = 14287962
$ws.Cells.Item(138, 10).Value = 3704.0715
$ws.Cells.Item(138, 11).Value = 42863886
$ws.Cells.Item(138, 12).Value = 11112.2145
$ws.Cells.Item(138, 13).Value = -42858746
$ws.Cells.Item(138, 14).Value = -21392.2145
$ws.Cells.Item(141, 8).Value = 6000
$ws.Cells.Item(141, 9).Value = 6000
$ws.Cells.Item(141, 10).Value = 0
$ws.Cells.Item(141, 11).Value = 18000
$ws.Cells.Item(141, 12).Value = 0
$ws.Cells.Item(141, 13).Value = -12820
$ws.Cells.Item(141, 14).ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 1831.4
$ws.Cells.Item(32, 9).Value = 1831.4
$ws.Cells.Item(32, 11).Value = 1831.4
$ws.Cells.Item(32, 13).Value = -1544.4
$ws.Cells.Item(92, 8).Value = 34999.5
$ws.Cells.Item(92, 10).Value = 34999.5
$ws.Cells.Item(92, 12).Value = 34999.5
$ws.Cells.Item(92, 14).Value = -39991.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 171.33333
$ws.Cells.Item(22, 9).Value = 171.33333
$ws.Cells.Item(22, 11).Value = 171.33333
$ws.Cells.Item(22, 13).Value = 1.666670000000011

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 7662.2
$ws.Cells.Item(31, 9).Value = 5827.75
$ws.Cells.Item(31, 10).Value = 15000
$ws.Cells.Item(31, 11).Value = 5827.75
$ws.Cells.Item(31, 12).Value = 15000
$ws.Cells.Item(31, 13).Value = -5532.75
$ws.Cells.Item(31, 14).Value = -15590
$ws.Cells.Item(34, 8).Value = 7662.2
$ws.Cells.Item(34, 9).Value = 5827.75
$ws.Cells.Item(34, 10).Value = 15000
$ws.Cells.Item(34, 11).Value = 5827.75
$ws.Cells.Item(34, 12).Value = 15000
$ws.Cells.Item(34, 13).Value = -5625.75
$ws.Cells.Item(34, 14).Value = -15404
$ws.Cells.Item(74, 8).Value = 28997.5
$ws.Cells.Item(74, 10).Value = 28997.5
$ws.Cells.Item(74, 12).Value = 28997.5
$ws.Cells.Item(74, 14).Value = -30745.5
$ws.Cells.Item(77, 8).Value = 28997.5
$ws.Cells.Item(77, 10).Value = 28997.5
$ws.Cells.Item(77, 12).Value = 86992.5
$ws.Cells.Item(77, 14).Value = -95728.5
$ws.Cells.Item(95, 8).Value = 22000
$ws.Cells.Item(95, 10).Value = 22000
$ws.Cells.Item(95, 12).Value = 22000
$ws.Cells.Item(95, 14).Value = -27492
$ws.Cells.Item(96, 8).Value = 15000
$ws.Cells.Item(96, 10).Value = 15000
$ws.Cells.Item(96, 12).Value = 15000
$ws.Cells.Item(96, 14).Value = -20492
$ws.Cells.Item(141, 8).Value = 501554.84
$ws.Cells.Item(141, 10).Value = 501554.84
$ws.Cells.Item(141, 12).Value = 501554.84
$ws.Cells.Item(141, 14).Value = -511914.84

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(11, 8).Value = 2912.1428
$ws.Cells.Item(11, 9).Value = 561.6667
$ws.Cells.Item(11, 10).Value = 4675
$ws.Cells.Item(11, 11).Value = 1685.0001
$ws.Cells.Item(11, 12).Value = 14025
$ws.Cells.Item(11, 13).Value = -1545.0001
$ws.Cells.Item(11, 14).Value = -14305
$ws.Cells.Item(97, 8).Value = 863.3333
$ws.Cells.Item(97, 9).Value = 545
$ws.Cells.Item(97, 10).Value = 1500
$ws.Cells.Item(97, 11).Value = 1635
$ws.Cells.Item(97, 12).Value = 4500
$ws.Cells.Item(97, 13).Value = -1139
$ws.Cells.Item(97, 14).Value = -5492
$ws.Cells.Item(122, 8).Value = 900
$ws.Cells.Item(122, 10).Value = 650
$ws.Cells.Item(122, 12).Value = 5850
$ws.Cells.Item(122, 14).Value = -10750
$ws.Cells.Item(131, 8).Value = 8832.833000000001
$ws.Cells.Item(131, 10).Value = 6666
$ws.Cells.Item(131, 12).Value = 19998
$ws.Cells.Item(131, 14).Value = -30078
$ws.Cells.Item(132, 8).Value = 350
$ws.Cells.Item(132, 9).Value = 350
$ws.Cells.Item(132, 11).Value = 3150
$ws.Cells.Item(132, 13).Value = -620

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(7, 8).Value = 9900
$ws.Cells.Item(7, 10).Value = 9900
$ws.Cells.Item(7, 12).Value = 9900
$ws.Cells.Item(7, 14).Value = -10124
$ws.Cells.Item(8, 8).Value = 9900
$ws.Cells.Item(8, 10).Value = 9900
$ws.Cells.Item(8, 12).Value = 9900
$ws.Cells.Item(8, 14).Value = -10178
$ws.Cells.Item(132, 8).Value = 1607.0834
$ws.Cells.Item(132, 9).Value = 1607.0834
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 4821.2502
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).Value = -2291.2502
$ws.Cells.Item(132, 14).ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(3, 8).Value = 0
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 12).Value = 0
$ws.Cells.Item(3, 14).ClearContents()
$ws.Cells.Item(15, 8).Value = 0
$ws.Cells.Item(15, 10).Value = 0
$ws.Cells.Item(15, 12).Value = 0
$ws.Cells.Item(15, 14).ClearContents()
$ws.Cells.Item(20, 8).Value = 16600
$ws.Cells.Item(20, 10).Value = 16600
$ws.Cells.Item(20, 12).Value = 16600
$ws.Cells.Item(20, 14).Value = -17052
$ws.Cells.Item(21, 8).Value = 18689
$ws.Cells.Item(21, 10).Value = 18689
$ws.Cells.Item(21, 12).Value = 18689
$ws.Cells.Item(21, 14).Value = -19037
$ws.Cells.Item(24, 8).Value = 19900
$ws.Cells.Item(24, 10).Value = 19900
$ws.Cells.Item(24, 12).Value = 19900
$ws.Cells.Item(24, 14).Value = -20586
$ws.Cells.Item(46, 8).Value = 1830.2667
$ws.Cells.Item(46, 9).Value = 1868.875
$ws.Cells.Item(46, 10).Value = 1786.1428
$ws.Cells.Item(46, 11).Value = 1868.875
$ws.Cells.Item(46, 12).Value = 1786.1428
$ws.Cells.Item(46, 13).Value = -1680.875
$ws.Cells.Item(46, 14).Value = -2162.1428

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(21, 8).Value = 22800
$ws.Cells.Item(21, 10).Value = 22800
$ws.Cells.Item(21, 12).Value = 22800
$ws.Cells.Item(21, 14).Value = -23270
$ws.Cells.Item(35, 8).Value = 22800
$ws.Cells.Item(35, 10).Value = 22800
$ws.Cells.Item(35, 12).Value = 22800
$ws.Cells.Item(35, 14).Value = -23380
$ws.Cells.Item(63, 8).Value = 24554.5
$ws.Cells.Item(63, 10).Value = 24554.5
$ws.Cells.Item(63, 12).Value = 24554.5
$ws.Cells.Item(63, 14).Value = -25802.5
$ws.Cells.Item(66, 8).Value = 24554.5
$ws.Cells.Item(66, 10).Value = 24554.5
$ws.Cells.Item(66, 12).Value = 73663.5
$ws.Cells.Item(66, 14).Value = -79903.5
$ws.Cells.Item(69, 8).Value = 27100
$ws.Cells.Item(69, 10).Value = 27100
$ws.Cells.Item(69, 12).Value = 27100
$ws.Cells.Item(69, 14).Value = -28598
$ws.Cells.Item(72, 8).Value = 27100
$ws.Cells.Item(72, 10).Value = 27100
$ws.Cells.Item(72, 12).Value = 81300
$ws.Cells.Item(72, 14).Value = -88788
$ws.Cells.Item(107, 8).Value = 1134.1052
$ws.Cells.Item(107, 9).Value = 1168.9
$ws.Cells.Item(107, 10).Value = 1095.4445
$ws.Cells.Item(107, 11).Value = 3506.7
$ws.Cells.Item(107, 12).Value = 3286.3335
$ws.Cells.Item(107, 13).Value = -1586.7
$ws.Cells.Item(107, 14).Value = -7126.333500000001
$ws.Cells.Item(109, 8).Value = 45000
$ws.Cells.Item(109, 9).Value = 45000
$ws.Cells.Item(109, 11).Value = 45000
$ws.Cells.Item(109, 13).Value = -43613
$ws.Cells.Item(113, 8).Value = 1960
$ws.Cells.Item(113, 9).Value = 1467.4
$ws.Cells.Item(113, 11).Value = 4402.200000000001
$ws.Cells.Item(113, 13).Value = -2232.200000000001
$ws.Cells.Item(119, 8).Value = 0
$ws.Cells.Item(119, 10).Value = 0
$ws.Cells.Item(119, 12).Value = 0
$ws.Cells.Item(119, 14).ClearContents()
$ws.Cells.Item(124, 8).Value = 0
$ws.Cells.Item(124, 10).Value = 0
$ws.Cells.Item(124, 12).Value = 0
$ws.Cells.Item(124, 14).ClearContents()

